$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0

$ws.Range("H76").Value = 3471.818
$ws.Range("I76").Value = 3479
$ws.Range("J76").Value = 3400
$ws.Range("K76").Value = 3479
$ws.Range("L76").Value = 3400
$ws.Range("M76").Value = -3164

$ws.Range("H79").Value = 3471.818
$ws.Range("I79").Value = 3479
$ws.Range("J79").Value = 3400
$ws.Range("K79").Value = 3479
$ws.Range("L79").Value = 3400
$ws.Range("M79").Value = -2387

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H112").Value = 50001588
$ws.Range("I112").Value = 250000380
$ws.Range("J112").Value = 1890.8125
$ws.Range("K112").Value = 750001140
$ws.Range("L112").Value = 5672.4375
$ws.Range("M112").Value = -750000032
$ws.Range("N112").Value = -7888.4375

$ws.Range("H129").Value = 892.675
$ws.Range("I129").Value = 540.3
$ws.Range("J129").Value = 1010.13336
$ws.Range("K129").Value = 1620.9
$ws.Range("L129").Value = 3030.40008
$ws.Range("M129").Value = 3379.1
$ws.Range("N129").Value = -13030.40008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 2241.1428
$ws.Range("I46").Value = 1949.5
$ws.Range("J46").Value = 2630
$ws.Range("K46").Value = 1949.5
$ws.Range("L46").Value = 2630
$ws.Range("M46").Value = -1630.5
$ws.Range("N46").Value = -3268

$ws.Range("H63").Value = 2669.9167
$ws.Range("I63").Value = 2503.9
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 2503.9
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -1817.9
$ws.Range("N63").Value = -4872

$ws.Range("H66").Value = 2669.9167
$ws.Range("I66").Value = 2503.9
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 12519.5
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -9087.5
$ws.Range("N66").Value = -24364

$ws.Range("H97").Value = 2501008.8
$ws.Range("I97").Value = 2977171
$ws.Range("J97").Value = 1157.25
$ws.Range("K97").Value = 2977171
$ws.Range("L97").Value = 1157.25
$ws.Range("M97").Value = -2976675
$ws.Range("N97").Value = -2149.25

$ws.Range("H110").Value = 2002460
$ws.Range("I110").Value = 2500825
$ws.Range("J110").Value = 9000
$ws.Range("K110").Value = 2500825
$ws.Range("L110").Value = 9000
$ws.Range("M110").Value = -2498780
$ws.Range("N110").Value = -13090

$ws.Range("H122").Value = 7410064
$ws.Range("I122").Value = 2790.8333
$ws.Range("J122").Value = 12348246
$ws.Range("K122").Value = 8372.499899999999
$ws.Range("L122").Value = 37044738
$ws.Range("M122").Value = -5922.499899999999
$ws.Range("N122").Value = -37049638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 500010000
$ws.Range("I26").Value = 500010000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 500010000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -500009708

$ws.Range("H94").Value = 990.8182
$ws.Range("I94").Value = 1014.1429
$ws.Range("J94").Value = 950
$ws.Range("K94").Value = 1014.1429
$ws.Range("L94").Value = 950
$ws.Range("M94").Value = -563.1429000000001
$ws.Range("N94").Value = -1852

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 20000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 20000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 862.25
$ws.Range("I86").Value = 550
$ws.Range("J86").Value = 966.3333
$ws.Range("K86").Value = 1650
$ws.Range("L86").Value = 2898.9999
$ws.Range("M86").Value = -464
$ws.Range("N86").Value = -5270.9999

$ws.Range("H89").Value = 862.25
$ws.Range("I89").Value = 550
$ws.Range("J89").Value = 966.3333
$ws.Range("K89").Value = 4950
$ws.Range("L89").Value = 8696.9997
$ws.Range("M89").Value = 978
$ws.Range("N89").Value = -20552.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 197.6
$ws.Range("I2").Value = 197.6
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 197.6
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -84.59999999999999
$ws.Range("N2").ClearContents()

$ws.Range("H80").Value = 3948.1177
$ws.Range("I80").Value = 3360
$ws.Range("J80").Value = 4074.1428
$ws.Range("K80").Value = 3360
$ws.Range("L80").Value = 4074.1428
$ws.Range("M80").Value = -2362
$ws.Range("N80").Value = -6070.1428

$ws.Range("H83").Value = 3948.1177
$ws.Range("I83").Value = 3360
$ws.Range("J83").Value = 4074.1428
$ws.Range("K83").Value = 16800
$ws.Range("L83").Value = 20370.714
$ws.Range("M83").Value = -11808
$ws.Range("N83").Value = -30354.714

$ws.Range("H97").Value = 1298.3334
$ws.Range("I97").Value = 1521.1111
$ws.Range("J97").Value = 630
$ws.Range("K97").Value = 1521.1111
$ws.Range("L97").Value = 630
$ws.Range("M97").Value = -1025.1111
$ws.Range("N97").Value = -1622

$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H122").Value = 4201.75
$ws.Range("I122").Value = 3101
$ws.Range("J122").Value = 7504
$ws.Range("K122").Value = 9303
$ws.Range("L122").Value = 22512
$ws.Range("M122").Value = -6853

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 400
$ws.Range("I11").Value = 400
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 400
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -260

$ws.Range("H40").Value = 2500
$ws.Range("I40").Value = 2466.6667
$ws.Range("J40").Value = 2700
$ws.Range("K40").Value = 2466.6667
$ws.Range("L40").Value = 2700
$ws.Range("M40").Value = -2330.6667
$ws.Range("N40").Value = -2972

$ws.Range("H68").Value = 836.92
$ws.Range("I68").Value = 706.6
$ws.Range("J68").Value = 1358.2
$ws.Range("K68").Value = 706.6
$ws.Range("L68").Value = 1358.2
$ws.Range("M68").Value = 42.39999999999998
$ws.Range("N68").Value = -2856.2

$ws.Range("H71").Value = 836.92
$ws.Range("I71").Value = 706.6
$ws.Range("J71").Value = 1358.2
$ws.Range("K71").Value = 3533
$ws.Range("L71").Value = 6791
$ws.Range("M71").Value = 211
$ws.Range("N71").Value = -14279

$ws.Range("H82").Value = 5015.364
$ws.Range("I82").Value = 999.6667
$ws.Range("J82").Value = 6521.25
$ws.Range("K82").Value = 999.6667
$ws.Range("L82").Value = 6521.25
$ws.Range("M82").Value = -638.6667
$ws.Range("N82").Value = -7243.25

$ws.Range("H85").Value = 5015.364
$ws.Range("I85").Value = 999.6667
$ws.Range("J85").Value = 6521.25
$ws.Range("K85").Value = 999.6667
$ws.Range("L85").Value = 6521.25
$ws.Range("M85").Value = 248.3333
$ws.Range("N85").Value = -9017.25

$ws.Range("H100").Value = 51403
$ws.Range("I100").Value = 78082.30499999999
$ws.Range("J100").Value = 1855.7142
$ws.Range("K100").Value = 78082.30499999999
$ws.Range("L100").Value = 1855.7142
$ws.Range("M100").Value = -77541.30499999999
$ws.Range("N100").Value = -2937.7142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3402.0557
$ws.Range("I62").Value = 3013.7646
$ws.Range("J62").Value = 10003
$ws.Range("K62").Value = 3013.7646
$ws.Range("L62").Value = 10003
$ws.Range("M62").Value = -2389.7646
$ws.Range("N62").Value = -11251

$ws.Range("H65").Value = 3402.0557
$ws.Range("I65").Value = 3013.7646
$ws.Range("J65").Value = 10003
$ws.Range("K65").Value = 15068.823
$ws.Range("L65").Value = 50015
$ws.Range("M65").Value = -11948.823
$ws.Range("N65").Value = -56255

$ws.Range("H81").Value = 2374.0527
$ws.Range("I81").Value = 1300
$ws.Range("J81").Value = 2757.6428
$ws.Range("K81").Value = 2600
$ws.Range("L81").Value = 5515.2856
$ws.Range("M81").Value = -1539

$ws.Range("H84").Value = 2374.0527
$ws.Range("I84").Value = 1300
$ws.Range("J84").Value = 2757.6428
$ws.Range("K84").Value = 13000
$ws.Range("L84").Value = 27576.428
$ws.Range("M84").Value = -7696

$ws.Range("H96").Value = 9611
$ws.Range("I96").Value = 6199.8
$ws.Range("J96").Value = 13875
$ws.Range("K96").Value = 6199.8
$ws.Range("L96").Value = 13875
$ws.Range("M96").Value = -4826.8
$ws.Range("N96").Value = -16621

$ws.Range("H109").Value = 24792.334
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 24792.334
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 24792.334
$ws.Range("N109").Value = -27566.334

$ws.Range("H122").Value = 2073.8965
$ws.Range("I122").Value = 1564.3334
$ws.Range("J122").Value = 4519.8
$ws.Range("K122").Value = 4693.0002
$ws.Range("L122").Value = 13559.4
$ws.Range("M122").Value = -2243.0002
$ws.Range("N122").Value = -18459.4
